$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap data between row 5 and row 6 (columns A, Q, R) ---
$a5 = $ws.Range("A5").Value2
$q5 = $ws.Range("Q5").Value2
$r5 = $ws.Range("R5").Value2

$a6 = $ws.Range("A6").Value2
$q6 = $ws.Range("Q6").Value2
$r6 = $ws.Range("R6").Value2

$ws.Range("A5").Value = $a6
$ws.Range("Q5").Value = $q6
$ws.Range("R5").Value = $r6

$ws.Range("A6").Value = $a5
$ws.Range("Q6").Value = $q5
$ws.Range("R6").Value = $r5

# --- Swap data between row 9 and row 10 (columns A, B, E, F, G, H, Q, R, S) ---
$a9 = $ws.Range("A9").Value2
$b9 = $ws.Range("B9").Value2
$e9 = $ws.Range("E9").Value2
$f9 = $ws.Range("F9").Value2
$g9 = $ws.Range("G9").Value2
$h9 = $ws.Range("H9").Value2
$q9 = $ws.Range("Q9").Value2
$r9 = $ws.Range("R9").Value2
$s9 = $ws.Range("S9").Value2

$a10 = $ws.Range("A10").Value2
$b10 = $ws.Range("B10").Value2
$e10 = $ws.Range("E10").Value2
$f10 = $ws.Range("F10").Value2
$g10 = $ws.Range("G10").Value2
$h10 = $ws.Range("H10").Value2
$q10 = $ws.Range("Q10").Value2
$r10 = $ws.Range("R10").Value2
$s10 = $ws.Range("S10").Value2

$ws.Range("A9").Value = $a10
$ws.Range("B9").Value = $b10
$ws.Range("E9").Value = $e10
$ws.Range("F9").Value = $f10
$ws.Range("G9").Value = $g10
$ws.Range("H9").Value = $h10
$ws.Range("Q9").Value = $q10
$ws.Range("R9").Value = $r10
$ws.Range("S9").Value = $s10

$ws.Range("A10").Value = $a9
$ws.Range("B10").Value = $b9
$ws.Range("E10").Value = $e9
$ws.Range("F10").Value = $f9
$ws.Range("G10").Value = $g9
$ws.Range("H10").Value = $h9
$ws.Range("Q10").Value = $q9
$ws.Range("R10").Value = $r9
$ws.Range("S10").Value = $s9

# --- Swap data between row 18 and row 19 (columns A, P, Q, R, S) ---
$a18 = $ws.Range("A18").Value2
$p18 = $ws.Range("P18").Value2
$q18 = $ws.Range("Q18").Value2
$r18 = $ws.Range("R18").Value2
$s18 = $ws.Range("S18").Value2

$a19 = $ws.Range("A19").Value2
$p19 = $ws.Range("P19").Value2
$q19 = $ws.Range("Q19").Value2
$r19 = $ws.Range("R19").Value2
$s19 = $ws.Range("S19").Value2

$ws.Range("A18").Value = $a19
$ws.Range("P18").Value = $p19
$ws.Range("Q18").Value = $q19
$ws.Range("R18").Value = $r19
$ws.Range("S18").Value = $s19

$ws.Range("A19").Value = $a18
$ws.Range("P19").Value = $p18
$ws.Range("Q19").Value = $q18
$ws.Range("R19").Value = $r18
$ws.Range("S19").Value = $s18
